$d = $word.ActiveDocument

# 1) The last paragraph of the document currently has its text split into
#    two runs ("...share th" / "e data it produces. ") with a collapsed
#    "_GoBack" bookmark sitting at the split point. Re-typing the full
#    sentence (the Find/Replace range spans over the bookmark) merges the
#    runs back into one and removes that now-stale bookmark.
$d.Content.Find.Execute(
    "Data producer, data consumer. Can be subscribed to. Can subscribe to other observables. All subscribers to a subject share the data it produces. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Data producer, data consumer. Can be subscribed to. Can subscribe to other observables. All subscribers to a subject share the data it produces. ",
    2) | Out-Null

# 2) Append the new notes (two more glossary-style entries) at the very end
#    of the document, reproducing Word's own WordprocessingML (including the
#    spell-check proofErr markers Word leaves around the technical terms),
#    and put the "_GoBack" bookmark back at the new end of the document,
#    where Word leaves it after the last edit.
$endRng = $d.Range($d.Content.End - 1, $d.Content.End - 1)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr>
</w:p>
<w:p>
<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:b/></w:rPr><w:t>HttpUrlEncodingCodec</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">: </w:t></w:r>
<w:r><w:t>Encodes the URL strings. I overwrote this for our dotnet core REST API once.</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:b/></w:rPr></w:pPr>
</w:p>
<w:p>
<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:b/></w:rPr><w:t>HttpXsrfTokenExtractor</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r>
<w:r><w:t xml:space="preserve"> Used to extract an </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>xsrf</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> token so it can be combined into the next request.</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$endRng.InsertXML($xml)
